$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.597.82'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').Value = '3.453.49'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '584.20'
$ws.Range('E5').Value = '  -1.31%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '176.12'
$ws.Range('E6').Value = '  -2.24%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.599'
$ws.Range('E8').Value = '  -2.57%  '
$ws.Range('D9').Value = '3.455.32'
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.131'
$ws.Range('E10').Value = '  -6.45%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.88'
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.413'
$ws.Range('E12').Value = '  -4.10%  '
$ws.Range('D13').Value = '4.056.97'
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '30.19'
$ws.Range('E15').Value = '  -5.48%  '
$ws.Range('D16').Value = '66.590.42'
$ws.Range('E16').Value = '  -0.73%  '
$ws.Range('E17').Value = '  -3.13%  '
$ws.Range('D18').Value = '3.455.96'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.92'
$ws.Range('E19').Value = '  -4.70%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.79'
$ws.Range('E20').Value = '  -2.71%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '377.37'
$ws.Range('E21').Value = '  -3.87%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.79'
$ws.Range('E22').Value = '  -1.81%  '
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '72.24'
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.75'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.534'
$ws.Range('E26').Value = '  -0.49%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000119'
$ws.Range('E27').Value = '  -1.87%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.70'
$ws.Range('E28').Value = '  -5.82%  '
$ws.Range('E29').Value = '  +1.18%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '24.18'
$ws.Range('E31').Value = '  +3.05%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.82'
$ws.Range('E32').Value = '  -5.11%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.99'
$ws.Range('E33').Value = '  -2.75%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  -6.44%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '7.14'
$ws.Range('E36').Value = '  -2.50%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.56'
$ws.Range('E37').Value = '  -1.50%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '160.85'
$ws.Range('E38').Value = '  -1.65%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '29.28'
$ws.Range('E39').Value = '  +11.68%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.892'
$ws.Range('E40').Value = '  +1.92%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.65'
$ws.Range('E41').Value = '  -6.17%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.77'
$ws.Range('E42').Value = '  -5.17%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.45'
$ws.Range('E43').Value = '  -4.61%  '
$ws.Range('D44').Value = '2.723.22'
$ws.Range('E44').Value = '  -1.00%  '
$ws.Range('E45').Value = '  -6.16%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0689'
$ws.Range('E46').Value = '  -4.23%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '40.71'
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0293'
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '24.24'
$ws.Range('E49').Value = '  -7.42%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '308.90'
$ws.Range('E50').Value = '  -4.99%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.823'
$ws.Range('E51').Value = '  -2.04%  '
